$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Materials": fix the tooltip text and convert the weight/atom %
# composition values so that percentages are always stored as 0-1
# fractions (with a 0.00% number format) instead of raw 0-100 numbers.
# ---------------------------------------------------------------------
$materials = $wb.Worksheets.Item("Materials")

# Tooltip text correction
$materials.Range("B3").Value = "Atom & Weight"

# Apply the percentage number format to the whole composition block
# (Target/Min/Max for both Weight % and Atom %, rows 13-16)
$materials.Range("B13:G16").NumberFormat = "0.00%"

# Row 13 - Carbon (C)
$materials.Range("D13").Value = 0.2
$materials.Range("G13").ClearContents()
$materials.Range("E13").Value = 0.2

# Row 14 - Chromium (Cr)
$materials.Range("B14").Value = 0.05
$materials.Range("D14").Value = 0.05

# Row 15 - Tungsten (W)
$materials.Range("C15").Value = 0.05
$materials.Range("D15").Value = 0.15

# Row 16 - Iron (Fe)
$materials.Range("B16").Value = 0.65
$materials.Range("C16").Value = 0.65
$materials.Range("D16").Value = 1

# Reflect the last-selected cell on this sheet
$materials.Range("E14").Select()

# ---------------------------------------------------------------------
# Sheet "Operating Conditions": Availability Factor was stored as the
# text "66%" - store it as a real number (0.66) with a percent format.
# ---------------------------------------------------------------------
$opConditions = $wb.Worksheets.Item("Operating Conditions")
$opConditions.Range("B7").NumberFormat = "0.00%"
$opConditions.Range("B7").Value = 0.66

$opConditions.Range("B7").Select()
